$wb = $excel.ActiveWorkbook

# ---- Sheet1 ("Sheet1"): mark BNS2 (row 3) and BNS4 (row 4) as approved ----
$ws1 = $wb.Worksheets.Item(1)

# Row 3 -> SubID BNS2
$ws1.Range("H3").Value = $true
$ws1.Range("N3").Value = "Ban Điều hành đã duyệt"

# Row 4 -> SubID BNS4
$ws1.Range("H4").Value = $true
$ws1.Range("N4").Value = "Ban Điều hành đã duyệt"

# ---- Sheet2 ("Ban Nhân Sự"): renumber + replace row 5 + append row 6 ----
$ws2 = $wb.Worksheets.Item(2)

# Existing rows just get their running ID decremented by 1
$ws2.Range("A2").Value = 5
$ws2.Range("A3").Value = 4
$ws2.Range("A4").Value = 3

# Row 5 becomes the newly-approved BNS4 record (copied over from Sheet1),
# replacing the old "Quy trình đào tạo2" placeholder row
$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = "BNS4"
$ws2.Range("C5").Value = "Quy trình tuyển 5.3ưe344"
$ws2.Range("D5").Value = "Ban Nhân Sự"
$ws2.Range("E5").Value = "15/09/2022"
$ws2.Range("F5").Value = "<p>24</p>"
$ws2.Range("G5").Value = "https://www.plus2net.com44"
$ws2.Range("H5").Value = $false
$ws2.Range("I5").Value = "Chưa có phản hồi"
$ws2.Range("J5").Value = $false
$ws2.Range("K5").Value = "Chưa có phản hồi"
$ws2.Range("L5").Value = $false
$ws2.Range("M5").Value = "Chưa có phản hồi"
$ws2.Range("N5").Value = "Ban Điều hành đã duyệt"

# New row 6: the newly-approved BNS2 record (copied over from Sheet1)
$ws2.Range("A6").Value = 1
$ws2.Range("B6").Value = "BNS2"
$ws2.Range("C6").Value = "Ngô Xuân Hinh123"
$ws2.Range("D6").Value = "Ban Nhân Sự"
$ws2.Range("E6").Value = "15/09/2022"
$ws2.Range("F6").Value = "<p>123<br></p>"
$ws2.Range("G6").Value = "https://www.plus2net.com43334"
$ws2.Range("H6").Value = $false
$ws2.Range("I6").Value = "Chưa có phản hồi"
$ws2.Range("J6").Value = $false
$ws2.Range("K6").Value = "Chưa có phản hồi"
$ws2.Range("L6").Value = $false
$ws2.Range("M6").Value = "Chưa có phản hồi"
$ws2.Range("N6").Value = "Ban Điều hành đã duyệt"
